{"js": "// Replace the multiplication-problem text runs in the practice table with\n// the new values from the commit, matching each old string exactly once\n// (all the old/new values here are unique within the document, so a\n// straightforward ordered search/replace is safe and unambiguous).\nconst replacements = [\n  [\"593\u00d75=\", \"170\u00d79=\"],\n  [\"817\u00d72=\", \"864\u00d78=\"],\n  [\"232\u00d78=\", \"336\u00d72=\"],\n  [\"133\u00d73=\", \"530\u00d74=\"],\n  [\"146\u00d78=\", \"667\u00d73=\"],\n  [\"964\u00d75=\", \"185\u00d72=\"],\n  [\"879\u00d79=\", \"436\u00d78=\"],\n  [\"739\u00d73=\", \"467\u00d76=\"],\n  [\"436\u00d76=\", \"145\u00d79=\"],\n  [\"485\u00d79=\", \"320\u00d79=\"],\n  [\"171\u00d76=\", \"177\u00d77=\"],\n  [\"961\u00d77=\", \"279\u00d72=\"],\n  [\"443\u00d77=\", \"312\u00d78=\"],\n  [\"715\u00d73=\", \"849\u00d76=\"],\n  [\"826\u00d74=\", \"556\u00d72=\"],\n  [\"876\u00d75=\", \"699\u00d75=\"],\n  [\"119\u00d73=\", \"384\u00d75=\"],\n  [\"553\u00d78=\", \"987\u00d75=\"],\n  [\"128\u00d76=\", \"759\u00d72=\"],\n  [\"800\u00d79=\", \"599\u00d75=\"],\n  [\"288\u00d78=\", \"783\u00d78=\"],\n  [\"957\u00d75=\", \"421\u00d75=\"],\n  [\"465\u00d74=\", \"503\u00d74=\"],\n  [\"490\u00d72=\", \"712\u00d75=\"],\n  [\"591\u00d72=\", \"149\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-problem text runs in the practice table with\n# the new values from the commit. All old/new values are unique within the\n# document, so Find/Replace (ReplaceAll) for each pair is safe and\n# unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"593\u00d75=\", \"170\u00d79=\"),\n    @(\"817\u00d72=\", \"864\u00d78=\"),\n    @(\"232\u00d78=\", \"336\u00d72=\"),\n    @(\"133\u00d73=\", \"530\u00d74=\"),\n    @(\"146\u00d78=\", \"667\u00d73=\"),\n    @(\"964\u00d75=\", \"185\u00d72=\"),\n    @(\"879\u00d79=\", \"436\u00d78=\"),\n    @(\"739\u00d73=\", \"467\u00d76=\"),\n    @(\"436\u00d76=\", \"145\u00d79=\"),\n    @(\"485\u00d79=\", \"320\u00d79=\"),\n    @(\"171\u00d76=\", \"177\u00d77=\"),\n    @(\"961\u00d77=\", \"279\u00d72=\"),\n    @(\"443\u00d77=\", \"312\u00d78=\"),\n    @(\"715\u00d73=\", \"849\u00d76=\"),\n    @(\"826\u00d74=\", \"556\u00d72=\"),\n    @(\"876\u00d75=\", \"699\u00d75=\"),\n    @(\"119\u00d73=\", \"384\u00d75=\"),\n    @(\"553\u00d78=\", \"987\u00d75=\"),\n    @(\"128\u00d76=\", \"759\u00d72=\"),\n    @(\"800\u00d79=\", \"599\u00d75=\"),\n    @(\"288\u00d78=\", \"783\u00d78=\"),\n    @(\"957\u00d75=\", \"421\u00d75=\"),\n    @(\"465\u00d74=\", \"503\u00d74=\"),\n    @(\"490\u00d72=\", \"712\u00d75=\"),\n    @(\"591\u00d72=\", \"149\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
